# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (fund-holdings snapshot, same layout as the
# existing "2021-Q4" sheet) positioned between "2021-Q4" and "总计", and
# records the new quarter in the "总计" (totals) summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet with the fund-holdings detail, right before
#    the "总计" sheet.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"
$newSheet.Move($wb.Worksheets.Item("总计"))

$q1 = $wb.Worksheets.Item("2022-Q1")

# Same headers + row-2 layout (values, types, styles) as "2021-Q4" - copy it
# wholesale, then overwrite the figures that are different for this quarter.
$src.Range("B1:H2").Copy($q1.Range("B1"))
$src.Range("A2").Copy($q1.Range("A2"))

# Figures that changed for 2022-Q1 (kept as text, same as the source sheet).
$q1.Range("D2").Value = "'0.04"
$q1.Range("E2").Value = "'22.00"
$q1.Range("F2").Value = "'0.92"
$q1.Range("G2").Value = "'0.0004"

# The leading "'" marks those cells as text but also tags them with a
# quote-prefix style; strip that back off so they match the plain styling of
# the rest of the row (use a blank, never-touched cell as the format source).
$blank = $q1.Range("Z100")
$blank.Copy()
$q1.Range("D2:G2").PasteSpecial(-4122)
$blank.ClearContents()

# ---------------------------------------------------------------------------
# 2. Record the new quarter in the "总计" summary sheet: insert a row for
#    "2022-Q1" above the existing "2021-Q4" row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "'2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

# The old "2021-Q4" row shifted down to row 3 - bump its running index.
$total.Range("A3").Value = 1

# Clean up styling picked up from the row insert / quoted text entry so the
# new cells match the look of the existing data row.
$blank2 = $total.Range("Z100")
$blank2.Copy()
$total.Range("B2:D2").PasteSpecial(-4122)
$blank2.ClearContents()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
